# Update profit-tracking figures across the Titan_Profits leve sheets
# (scheduled price-refresh run: currentAveragePrice* / LevePrice* / LeveProfit* columns)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43: Growth Formula Gamma
$ws.Range("H43").Value = 923
$ws.Range("I43").Value = 866.6667
$ws.Range("J43").Value = 956.8
$ws.Range("K43").Value = 866.6667
$ws.Range("L43").Value = 956.8
$ws.Range("M43").Value = -797.6667
$ws.Range("N43").Value = -1094.8

# Row 51: Shark Oil
$ws.Range("H51").Value = 2179.6
$ws.Range("I51").Value = 1133
$ws.Range("J51").Value = 3749.5
$ws.Range("K51").Value = 1133
$ws.Range("L51").Value = 3749.5
$ws.Range("M51").Value = -649
$ws.Range("N51").Value = -4717.5

# Row 111: Grade 1 Dexterity Alkahest
$ws.Range("H111").Value = 1621.8
$ws.Range("I111").Value = 1209.6666
$ws.Range("J111").Value = 2240
$ws.Range("K111").Value = 3628.9998
$ws.Range("L111").Value = 6720
$ws.Range("M111").Value = -561.9998000000001
$ws.Range("N111").Value = -12854

# Row 129: Commanding Craftsman's Draught
$ws.Range("H129").Value = 1232.5883
$ws.Range("I129").Value = 546
$ws.Range("J129").Value = 1443.8462
$ws.Range("K129").Value = 1638
$ws.Range("L129").Value = 4331.5386
$ws.Range("M129").Value = 3362
$ws.Range("N129").Value = -14331.5386

# Row 131: Grade 5 Tincture of Mind
$ws.Range("H131").Value = 5995.7334
$ws.Range("I131").Value = 1405.6
$ws.Range("J131").Value = 15176
$ws.Range("K131").Value = 4216.799999999999
$ws.Range("L131").Value = 45528
$ws.Range("M131").Value = 823.2000000000007
$ws.Range("N131").Value = -55608

# Row 141: Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 1547.6316
$ws.Range("I141").Value = 1366
$ws.Range("J141").Value = 3666.6667
$ws.Range("K141").Value = 4098
$ws.Range("L141").Value = 11000.0001
$ws.Range("M141").Value = 1082
$ws.Range("N141").Value = -21360.0001

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Bronze Ingot
$ws.Range("H2").Value = 68355.47
$ws.Range("I2").Value = 78417.84
$ws.Range("J2").Value = 2950
$ws.Range("K2").Value = 78417.84
$ws.Range("L2").Value = 2950
$ws.Range("M2").Value = -78304.84
$ws.Range("N2").Value = -3176

# Row 61: Cobalt Ingot
$ws.Range("H61").Value = 1895.6129
$ws.Range("I61").Value = 1582.75
$ws.Range("J61").Value = 4815.6665
$ws.Range("K61").Value = 1582.75
$ws.Range("L61").Value = 4815.6665
$ws.Range("M61").Value = -1370.75
$ws.Range("N61").Value = -5239.6665

# Row 97: High Steel Ingot
$ws.Range("H97").Value = 6776.875
$ws.Range("I97").Value = 8164.615
$ws.Range("J97").Value = 763.3333
$ws.Range("K97").Value = 8164.615
$ws.Range("L97").Value = 763.3333
$ws.Range("M97").Value = -7668.615
$ws.Range("N97").Value = -1755.3333

# Row 102: Tama-hagane Ingot
$ws.Range("H102").Value = 1828
$ws.Range("I102").Value = 1570
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1570
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 52
$ws.Range("N102").Value = -5244

# Row 116: Titanbronze Ingot
$ws.Range("H116").Value = 68355.47
$ws.Range("I116").Value = 78417.84
$ws.Range("J116").Value = 2950
$ws.Range("K116").Value = 78417.84
$ws.Range("L116").Value = 2950
$ws.Range("M116").Value = -76123.84
$ws.Range("N116").Value = -7538

# Row 132: Mountain Chromite Ingot
$ws.Range("H132").Value = 2226.9722
$ws.Range("I132").Value = 1964.3704
$ws.Range("J132").Value = 3014.7778
$ws.Range("K132").Value = 5893.1112
$ws.Range("L132").Value = 9044.3334
$ws.Range("M132").Value = -3363.1112
$ws.Range("N132").Value = -14104.3334

# Row 136: Cobalt Tungsten Ingot
$ws.Range("H136").Value = 1895.6129
$ws.Range("I136").Value = 1582.75
$ws.Range("J136").Value = 4815.6665
$ws.Range("K136").Value = 4748.25
$ws.Range("L136").Value = 14446.9995
$ws.Range("M136").Value = -2198.25
$ws.Range("N136").Value = -19546.9995

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Bronze Ingot
$ws.Range("H3").Value = 68355.47
$ws.Range("I3").Value = 78417.84
$ws.Range("J3").Value = 2950
$ws.Range("K3").Value = 78417.84
$ws.Range("L3").Value = 2950
$ws.Range("M3").Value = -78303.84
$ws.Range("N3").Value = -3178

# Row 134: Ruthenium Ingot
$ws.Range("H134").Value = 3466.6287
$ws.Range("I134").Value = 2178.55
$ws.Range("J134").Value = 5184.067
$ws.Range("K134").Value = 6535.650000000001
$ws.Range("L134").Value = 15552.201
$ws.Range("M134").Value = -4000.650000000001
$ws.Range("N134").Value = -20622.201

$ws = $wb.Worksheets.Item("CRP")
# Row 140: Claro Walnut Spear
$ws.Range("H140").Value = 56770
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 56770
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 56770
$ws.Range("N140").Value = -67130

$ws = $wb.Worksheets.Item("CUL")
# Row 107: Frantoio Oil
$ws.Range("H107").Value = 419.3
$ws.Range("I107").Value = 487.83334
$ws.Range("J107").Value = 389.92856
$ws.Range("K107").Value = 1463.50002
$ws.Range("L107").Value = 1169.78568
$ws.Range("M107").Value = 456.4999800000001
$ws.Range("N107").Value = -5009.78568

# Row 127: Carrot Nibbles
$ws.Range("H127").Value = 1963.3
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 1963.3
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 5889.9
$ws.Range("N127").Value = -15809.9

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Mythrite Ingot
$ws.Range("H70").Value = 9425
$ws.Range("I70").Value = 12450
$ws.Range("J70").Value = 6400
$ws.Range("K70").Value = 12450
$ws.Range("L70").Value = 6400
$ws.Range("M70").Value = -12180
$ws.Range("N70").Value = -6940

# Row 73: Mythrite Ingot
$ws.Range("H73").Value = 9425
$ws.Range("I73").Value = 12450
$ws.Range("J73").Value = 6400
$ws.Range("K73").Value = 12450
$ws.Range("L73").Value = 6400
$ws.Range("M73").Value = -11514
$ws.Range("N73").Value = -8272

# Row 102: Durium Ingot
$ws.Range("H102").Value = 3667.7144
$ws.Range("I102").Value = 3791.875
$ws.Range("J102").Value = 3502.1667
$ws.Range("K102").Value = 3791.875
$ws.Range("L102").Value = 3502.1667
$ws.Range("M102").Value = -2169.875
$ws.Range("N102").Value = -6746.1667

# Row 113: Manasilver Nugget
$ws.Range("H113").Value = 1842.8334
$ws.Range("I113").Value = 1842.5454
$ws.Range("J113").Value = 1843.2858
$ws.Range("K113").Value = 1842.5454
$ws.Range("L113").Value = 1843.2858
$ws.Range("M113").Value = 327.4546
$ws.Range("N113").Value = -6183.2858

# Row 132: Lar Ingot
$ws.Range("H132").Value = 2857.4
$ws.Range("I132").Value = 2423.8635
$ws.Range("J132").Value = 4591.5454
$ws.Range("K132").Value = 7271.5905
$ws.Range("L132").Value = 13774.6362
$ws.Range("M132").Value = -4741.5905
$ws.Range("N132").Value = -18834.6362

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Leather
$ws.Range("H7").Value = 3083.64
$ws.Range("I7").Value = 2000.2
$ws.Range("J7").Value = 3354.5
$ws.Range("K7").Value = 2000.2
$ws.Range("L7").Value = 3354.5
$ws.Range("M7").Value = -1888.2
$ws.Range("N7").Value = -3578.5

# Row 40: Toad Leather
$ws.Range("H40").Value = 3862.5
$ws.Range("I40").Value = 2314.2856
$ws.Range("J40").Value = 4500
$ws.Range("K40").Value = 2314.2856
$ws.Range("L40").Value = 4500
$ws.Range("M40").Value = -2178.2856
$ws.Range("N40").Value = -4772

# Row 126: Saiga Leather
$ws.Range("H126").Value = 3083.64
$ws.Range("I126").Value = 2000.2
$ws.Range("J126").Value = 3354.5
$ws.Range("K126").Value = 6000.6
$ws.Range("L126").Value = 10063.5
$ws.Range("M126").Value = -3530.6
$ws.Range("N126").Value = -15003.5

# Row 137: Br'aaxskin Halfgloves of Crafting
$ws.Range("H137").Value = 49143
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 49143
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 49143
$ws.Range("N137").Value = -59343

$ws = $wb.Worksheets.Item("WVR")
# Row 107: Bright Linen Yarn
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 873.35297
$ws.Range("J107").Value = 1215.3
$ws.Range("K107").Value = 2620.05891
$ws.Range("L107").Value = 3645.9
$ws.Range("M107").Value = -700.0589100000002
$ws.Range("N107").Value = -7485.9

# Row 108: Brightlinen Bottoms of Striking
$ws.Range("H108").Value = 48600
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 48600
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 48600
$ws.Range("N108").Value = -56280

# Row 132: Snow Cotton Cloth
$ws.Range("H132").Value = 11629786
$ws.Range("I132").Value = 16668331
$ws.Range("J132").Value = 2374.6155
$ws.Range("K132").Value = 50004993
$ws.Range("L132").Value = 7123.8465
$ws.Range("M132").Value = -50002463
$ws.Range("N132").Value = -12183.8465

# Row 136: Sarcenet Cloth
$ws.Range("H136").Value = 11941910
$ws.Range("I136").Value = 12860136
$ws.Range("J136").Value = 4977
$ws.Range("K136").Value = 38580408
$ws.Range("L136").Value = 14931
$ws.Range("M136").Value = -38577858
$ws.Range("N136").Value = -20031

